$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 64.961031
$ws.Range("H2").Value = 194.883093
$ws.Range("I2").Value = 0.7420823237482214
$ws.Range("J2").Value = 0.7420823237482211
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 64.961031
$ws.Range("N2").Value = 194.883093
$ws.Range("O2").Value = 0.7420823237482214
$ws.Range("P2").Value = 0.7420823237482211
$ws.Range("Q2").Value = 4219.935548582962
$ws.Range("R2").Value = 37979.41993724666
$ws.Range("S2").Value = 0.55068617521956
$ws.Range("T2").Value = 0.5506861752195596

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 64.961031
$ws.Range("H3").Value = 194.883093
$ws.Range("I3").Value = 0.7420823237482214
$ws.Range("J3").Value = 0.7420823237482211
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.342779
$ws.Range("N3").Value = 58.02833699999999
$ws.Range("O3").Value = 0.2209622317735119
$ws.Range("P3").Value = 0.2209622317735119
$ws.Range("Q3").Value = 1256.526866245149
$ws.Range("R3").Value = 11308.74179620634
$ws.Range("S3").Value = 0.1639721664150808
$ws.Range("T3").Value = 0.1639721664150807

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 64.961031
$ws.Range("H4").Value = 194.883093
$ws.Range("I4").Value = 0.7420823237482214
$ws.Range("J4").Value = 0.7420823237482211
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.168388
$ws.Range("N4").Value = 0.5051640000000001
$ws.Range("O4").Value = 0.001923580282020393
$ws.Range("P4").Value = 0.001923580282020393
$ws.Range("Q4").Value = 10.938658088028
$ws.Range("R4").Value = 98.44792279225203
$ws.Range("S4").Value = 0.001427454925597952
$ws.Range("T4").Value = 0.001427454925597951

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 64.961031
$ws.Range("H5").Value = 194.883093
$ws.Range("I5").Value = 0.7420823237482214
$ws.Range("J5").Value = 0.7420823237482211
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.066649
$ws.Range("N5").Value = 9.199947
$ws.Range("O5").Value = 0.0350318641962465
$ws.Range("P5").Value = 0.03503186419624649
$ws.Range("Q5").Value = 199.212680755119
$ws.Range("R5").Value = 1792.914126796071
$ws.Range("S5").Value = 0.02599652718798272
$ws.Range("T5").Value = 0.0259965271879827

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.342779
$ws.Range("H6").Value = 58.02833699999999
$ws.Range("I6").Value = 0.2209622317735119
$ws.Range("J6").Value = 0.2209622317735119
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 64.961031
$ws.Range("N6").Value = 194.883093
$ws.Range("O6").Value = 0.7420823237482214
$ws.Range("P6").Value = 0.7420823237482211
$ws.Range("Q6").Value = 1256.526866245149
$ws.Range("R6").Value = 11308.74179620634
$ws.Range("S6").Value = 0.1639721664150808
$ws.Range("T6").Value = 0.1639721664150807

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.342779
$ws.Range("H7").Value = 58.02833699999999
$ws.Range("I7").Value = 0.2209622317735119
$ws.Range("J7").Value = 0.2209622317735119
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.342779
$ws.Range("N7").Value = 58.02833699999999
$ws.Range("O7").Value = 0.2209622317735119
$ws.Range("P7").Value = 0.2209622317735119
$ws.Range("Q7").Value = 374.1430994428409
$ws.Range("R7").Value = 3367.287894985568
$ws.Range("S7").Value = 0.04882430787033121
$ws.Range("T7").Value = 0.04882430787033118

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.342779
$ws.Range("H8").Value = 58.02833699999999
$ws.Range("I8").Value = 0.2209622317735119
$ws.Range("J8").Value = 0.2209622317735119
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.168388
$ws.Range("N8").Value = 0.5051640000000001
$ws.Range("O8").Value = 0.001923580282020393
$ws.Range("P8").Value = 0.001923580282020393
$ws.Range("Q8").Value = 3.257091870252
$ws.Range("R8").Value = 29.313826832268
$ws.Range("S8").Value = 0.0004250385921107475
$ws.Range("T8").Value = 0.0004250385921107473

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.342779
$ws.Range("H9").Value = 58.02833699999999
$ws.Range("I9").Value = 0.2209622317735119
$ws.Range("J9").Value = 0.2209622317735119
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.066649
$ws.Range("N9").Value = 9.199947
$ws.Range("O9").Value = 0.0350318641962465
$ws.Range("P9").Value = 0.03503186419624649
$ws.Range("Q9").Value = 59.31751387757099
$ws.Range("R9").Value = 533.857624898139
$ws.Range("S9").Value = 0.007740718895989214
$ws.Range("T9").Value = 0.007740718895989209

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.168388
$ws.Range("H10").Value = 0.5051640000000001
$ws.Range("I10").Value = 0.001923580282020393
$ws.Range("J10").Value = 0.001923580282020393
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 64.961031
$ws.Range("N10").Value = 194.883093
$ws.Range("O10").Value = 0.7420823237482214
$ws.Range("P10").Value = 0.7420823237482211
$ws.Range("Q10").Value = 10.938658088028
$ws.Range("R10").Value = 98.44792279225203
$ws.Range("S10").Value = 0.001427454925597952
$ws.Range("T10").Value = 0.001427454925597951

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.168388
$ws.Range("H11").Value = 0.5051640000000001
$ws.Range("I11").Value = 0.001923580282020393
$ws.Range("J11").Value = 0.001923580282020393
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 19.342779
$ws.Range("N11").Value = 58.02833699999999
$ws.Range("O11").Value = 0.2209622317735119
$ws.Range("P11").Value = 0.2209622317735119
$ws.Range("Q11").Value = 3.257091870252
$ws.Range("R11").Value = 29.313826832268
$ws.Range("S11").Value = 0.0004250385921107475
$ws.Range("T11").Value = 0.0004250385921107473

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.168388
$ws.Range("H12").Value = 0.5051640000000001
$ws.Range("I12").Value = 0.001923580282020393
$ws.Range("J12").Value = 0.001923580282020393
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.168388
$ws.Range("N12").Value = 0.5051640000000001
$ws.Range("O12").Value = 0.001923580282020393
$ws.Range("P12").Value = 0.001923580282020393
$ws.Range("Q12").Value = 0.028354518544
$ws.Range("R12").Value = 0.2551906668960001
$ws.Range("S12").Value = 0.000003700161101377655
$ws.Range("T12").Value = 0.000003700161101377653

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.168388
$ws.Range("H13").Value = 0.5051640000000001
$ws.Range("I13").Value = 0.001923580282020393
$ws.Range("J13").Value = 0.001923580282020393
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.066649
$ws.Range("N13").Value = 9.199947
$ws.Range("O13").Value = 0.0350318641962465
$ws.Range("P13").Value = 0.03503186419624649
$ws.Range("Q13").Value = 0.516386891812
$ws.Range("R13").Value = 4.647482026308
$ws.Range("S13").Value = 0.00006738660321031595
$ws.Range("T13").Value = 0.00006738660321031591

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.066649
$ws.Range("H14").Value = 9.199947
$ws.Range("I14").Value = 0.0350318641962465
$ws.Range("J14").Value = 0.03503186419624649
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 64.961031
$ws.Range("N14").Value = 194.883093
$ws.Range("O14").Value = 0.7420823237482214
$ws.Range("P14").Value = 0.7420823237482211
$ws.Range("Q14").Value = 199.212680755119
$ws.Range("R14").Value = 1792.914126796071
$ws.Range("S14").Value = 0.02599652718798272
$ws.Range("T14").Value = 0.0259965271879827

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.066649
$ws.Range("H15").Value = 9.199947
$ws.Range("I15").Value = 0.0350318641962465
$ws.Range("J15").Value = 0.03503186419624649
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 19.342779
$ws.Range("N15").Value = 58.02833699999999
$ws.Range("O15").Value = 0.2209622317735119
$ws.Range("P15").Value = 0.2209622317735119
$ws.Range("Q15").Value = 59.31751387757099
$ws.Range("R15").Value = 533.857624898139
$ws.Range("S15").Value = 0.007740718895989214
$ws.Range("T15").Value = 0.007740718895989209

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.066649
$ws.Range("H16").Value = 9.199947
$ws.Range("I16").Value = 0.0350318641962465
$ws.Range("J16").Value = 0.03503186419624649
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.168388
$ws.Range("N16").Value = 0.5051640000000001
$ws.Range("O16").Value = 0.001923580282020393
$ws.Range("P16").Value = 0.001923580282020393
$ws.Range("Q16").Value = 0.516386891812
$ws.Range("R16").Value = 4.647482026308
$ws.Range("S16").Value = 0.00006738660321031595
$ws.Range("T16").Value = 0.00006738660321031591

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.066649
$ws.Range("H17").Value = 9.199947
$ws.Range("I17").Value = 0.0350318641962465
$ws.Range("J17").Value = 0.03503186419624649
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.066649
$ws.Range("N17").Value = 9.199947
$ws.Range("O17").Value = 0.0350318641962465
$ws.Range("P17").Value = 0.03503186419624649
$ws.Range("Q17").Value = 9.404336089201
$ws.Range("R17").Value = 84.639024802809
$ws.Range("S17").Value = 0.001227231509064258
$ws.Range("T17").Value = 0.001227231509064257
